# edit.ps1 - apply the "Helena" slot game review edits described by the diff.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title / meta-title text changes (two occurrences: the H1 heading and the
#    bold "SEO title" run near the end of the document). A global
#    Find/Replace cleanly swaps the <w:t> content while leaving every run's
#    formatting (and the leading empty run used throughout this document)
#    untouched.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Play Helena for Free - Slot Game Review", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Helena Slot Game for Free", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Meta description (the italic run at the very end of the document).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Explore the world of ancient Greek mythology with the beautifully designed Helena slot game. Place bets and win up to 40,000 coins. Play for free.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Read our review of Helena slot game and play for free. Enjoy stunning visuals and decent payouts.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Helpers for the "What we like" / "What we don't like" bullet list, which
# needs paragraphs inserted/merged/removed, not just re-texted. We rebuild
# affected paragraphs via Range.InsertXML so the familiar
#   <w:r/><w:r><w:t>...</w:t></w:r>
# shape (leading empty run + text run) used by every bullet in this document
# is preserved exactly.
# ---------------------------------------------------------------------------
function Get-ParagraphByText($doc, $text) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

# The exact <w:pPr> used by every "What we like" / "What we don't like" bullet.
$bulletPPr = '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>'

function Set-BulletParagraphText($para, $text) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $bulletPPr + '<w:r/><w:r><w:t>' + $text + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($xml) | Out-Null
}

# ---------------------------------------------------------------------------
# 3) "Autoplay option" -> "Autoplay option for convenience" and insert a new
#    bullet "Wide range of betting options" right after it.
# ---------------------------------------------------------------------------
$autoplayPara = Get-ParagraphByText $d "Autoplay option"
Set-BulletParagraphText $autoplayPara "Autoplay option for convenience"

$autoplayPara = Get-ParagraphByText $d "Autoplay option for convenience"
$autoplayPara.Range.InsertParagraphAfter() | Out-Null

# Locate the freshly-created (empty) ListBullet paragraph right after it and
# fill it in.
$found = $false
$newPara = $null
foreach ($p in $d.Paragraphs) {
    if ($found) {
        $newPara = $p
        break
    }
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Autoplay option for convenience") {
        $found = $true
    }
}
Set-BulletParagraphText $newPara "Wide range of betting options"

# ---------------------------------------------------------------------------
# 4) Merge "Decent payouts" + "High-value symbols" bullets into a single
#    "Decent payouts and high-value symbols" bullet (drop the second one).
# ---------------------------------------------------------------------------
$decentPara = Get-ParagraphByText $d "Decent payouts"
Set-BulletParagraphText $decentPara "Decent payouts and high-value symbols"

$highValuePara = Get-ParagraphByText $d "High-value symbols"
$highValuePara.Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 5) "What we don't like" bullet tweaks.
# ---------------------------------------------------------------------------
$rtpPara = Get-ParagraphByText $d "RTP is below average"
Set-BulletParagraphText $rtpPara "RTP slightly below average"

$volatilityPara = Get-ParagraphByText $d "Not ideal for players seeking high volatility"
Set-BulletParagraphText $volatilityPara "Limited number of paylines"

Write-Host "Edits applied."
